# cryptos.xlsx price/volume refresh
# Commit: Updated cryptos list on Fri Oct 13 16:40:57 UTC 2023 with GitHub Actions

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "26.942.25"
$ws.Range("E2").Value = "  +0.62%  "

# Row 3
$ws.Range("D3").Value = "1.554.09"
$ws.Range("E3").Value = "  +1.00%  "

# Row 4
$ws.Range("E4").Value = "  +0.45%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "207.11"
$ws.Range("E5").Value = "  +0.62%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.486"
$ws.Range("E6").Value = "  +0.99%  "

# Row 7
$ws.Range("E7").Value = "  +0.42%  "

# Row 8
$ws.Range("B8").Value = "Solana"
$ws.Range("C8").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "21.71"
$ws.Range("E8").Value = "  +2.11%  "

# Row 9
$ws.Range("B9").Value = "Cardano"
$ws.Range("C9").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.249"
$ws.Range("E9").Value = "  +1.82%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0589"
$ws.Range("E10").Value = "  +1.72%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0859"
$ws.Range("E11").Value = "  +0.66%  "

# Row 12
$ws.Range("D12").Value = "1.773.97"
$ws.Range("E12").Value = "  +0.89%  "

# Row 13
$ws.Range("D13").Value = "1.549.68"
$ws.Range("E13").Value = "  +0.72%  "

# Row 14
$ws.Range("E14").Value = "  +1.46%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.516"
$ws.Range("E15").Value = "  +1.80%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "61.87"
$ws.Range("E16").Value = "  +1.29%  "

# Row 17
$ws.Range("D17").Value = "26.936.49"
$ws.Range("E17").Value = "  +0.62%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "215.99"
$ws.Range("E18").Value = "  +1.79%  "

# Row 19
$ws.Range("D19").Value = "0.0₃0689"
$ws.Range("E19").Value = "  +0.14%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.23"
$ws.Range("E20").Value = "  +0.07%  "

# Row 21
$ws.Range("E21").Value = "  +0.44%  "

# Row 22
$ws.Range("E22").Value = "  +0.22%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.16"
$ws.Range("E23").Value = "  +1.50%  "

# Row 24
$ws.Range("E24").Value = "  +0.24%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "152.45"
$ws.Range("E25").Value = "  -0.23%  "

# Row 26
$ws.Range("E26").Value = "  +2.89%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "14.90"
$ws.Range("E27").Value = "  +0.44%  "

# Row 28
$ws.Range("E28").Value = "  +0.48%  "

# Row 29
$ws.Range("E29").Value = "  +1.26%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0464"
$ws.Range("E30").Value = "  +2.14%  "

# Row 31
$ws.Range("E31").Value = "  -0.46%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.22"
$ws.Range("E32").Value = "  -0.02%  "

# Row 33
$ws.Range("D33").Value = "1.417.13"
$ws.Range("E33").Value = "  +4.30%  "

# Row 34
$ws.Range("E34").Value = "  +3.02%  "

# Row 35
$ws.Range("E35").Value = "  +3.77%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.959"
$ws.Range("E36").Value = "  +3.26%  "

# Row 37
$ws.Range("E37").Value = "  +0.67%  "

# Row 38
$ws.Range("E38").Value = "  +1.06%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.521"
$ws.Range("E39").Value = "  -0.12%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.806"
$ws.Range("E40").Value = "  +1.12%  "

# Row 41
$ws.Range("E41").Value = "  +0.50%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.58"
$ws.Range("E42").Value = "  -2.70%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.990"
$ws.Range("E43").Value = "  -0.20%  "

# Row 44
$ws.Range("E44").Value = "  +3.60%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "63.68"
$ws.Range("E45").Value = "  +1.96%  "

# Row 46
$ws.Range("E46").Value = "  +0.69%  "

# Row 47
$ws.Range("D47").Value = "1.688.73"
$ws.Range("E47").Value = "  +0.64%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "86.16"
$ws.Range("E48").Value = "  +0.47%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0518"
$ws.Range("E49").Value = "  +1.52%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0957"
$ws.Range("E50").Value = "  +0.81%  "

# Row 51
$ws.Range("B51").Value = "BabyDogeCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D51").Value = "0.0₇0964"
$ws.Range("E51").Value = "  -0.97%  "

